$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so numeric-looking values
# (e.g. trailing zeros like "42.90") are preserved exactly, matching the
# original inline-string storage instead of being parsed into floats.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '30.094.38'
$ws.Range("E2").Value = '  +5.53%  '

# Row 3
$ws.Range("D3").Value = '1.924.66'
$ws.Range("E3").Value = '  +2.71%  '

# Row 4
$ws.Range("E4").Value = '  -0.87%  '

# Row 5
$ws.Range("D5").Value = '325.92'
$ws.Range("E5").Value = '  +3.10%  '

# Row 6
$ws.Range("E6").Value = '  -0.79%  '

# Row 7
$ws.Range("D7").Value = '0.5161'
$ws.Range("E7").Value = '  +1.59%  '

# Row 8
$ws.Range("D8").Value = '0.4005'
$ws.Range("E8").Value = '  +2.77%  '

# Row 9
$ws.Range("D9").Value = '0.08469'
$ws.Range("E9").Value = '  +1.15%  '

# Row 10
$ws.Range("D10").Value = '42.90'
$ws.Range("E10").Value = '  +2.62%  '

# Row 11
$ws.Range("D11").Value = '1.123'
$ws.Range("E11").Value = '  +1.77%  '

# Row 12
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '21.20'
$ws.Range("E12").Value = '  +3.82%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '6.337'
$ws.Range("E13").Value = '  +1.96%  '

# Row 14
$ws.Range("D14").Value = '1.918.19'
$ws.Range("E14").Value = '  +2.33%  '

# Row 15
$ws.Range("D15").Value = '7.345'
$ws.Range("E15").Value = '  +1.44%  '

# Row 16
$ws.Range("E16").Value = '  -0.91%  '

# Row 17
$ws.Range("D17").Value = '94.34'
$ws.Range("E17").Value = '  +3.44%  '

# Row 18
$ws.Range("E18").Value = '  +1.09%  '

# Row 19
$ws.Range("D19").Value = '0.06778'
$ws.Range("E19").Value = '  +1.04%  '

# Row 20
$ws.Range("D20").Value = '18.02'
$ws.Range("E20").Value = '  +1.76%  '

# Row 21
$ws.Range("E21").Value = '  -0.70%  '

# Row 22
$ws.Range("D22").Value = '6.064'
$ws.Range("E22").Value = '  +2.31%  '

# Row 23
$ws.Range("D23").Value = '30.092.02'
$ws.Range("E23").Value = '  +5.40%  '

# Row 24
$ws.Range("E24").Value = '  +1.06%  '

# Row 25
$ws.Range("D25").Value = '2.200'
$ws.Range("E25").Value = '  -1.63%  '

# Row 26
$ws.Range("D26").Value = '2.141.11'
$ws.Range("E26").Value = '  +2.58%  '

# Row 27
$ws.Range("D27").Value = '160.06'
$ws.Range("E27").Value = '  -0.98%  '

# Row 28
$ws.Range("E28").Value = '  +1.79%  '

# Row 29
$ws.Range("D29").Value = '2.469'
$ws.Range("E29").Value = '  +5.10%  '

# Row 30
$ws.Range("D30").Value = '128.87'
$ws.Range("E30").Value = '  +2.30%  '

# Row 31
$ws.Range("E31").Value = '  +3.31%  '

# Row 32
$ws.Range("D32").Value = '0.1059'
$ws.Range("E32").Value = '  +1.40%  '

# Row 33
$ws.Range("D33").Value = '6.080'
$ws.Range("E33").Value = '  +5.05%  '

# Row 34
$ws.Range("D34").Value = '3.648'
$ws.Range("E34").Value = '  +0.97%  '

# Row 35
$ws.Range("D35").Value = '0.02503'
$ws.Range("E35").Value = '  +1.83%  '

# Row 36
$ws.Range("D36").Value = '0.06607'
$ws.Range("E36").Value = '  +0.82%  '

# Row 37
$ws.Range("D37").Value = '0.2227'

# Row 38
$ws.Range("D38").Value = '1.244'
$ws.Range("E38").Value = '  +4.53%  '

# Row 39
$ws.Range("D39").Value = '9.014'
$ws.Range("E39").Value = '  +1.77%  '

# Row 40
$ws.Range("D40").Value = '5.204'
$ws.Range("E40").Value = '  +2.39%  '

# Row 41
$ws.Range("D41").Value = '0.6545'
$ws.Range("E41").Value = '  +1.89%  '

# Row 42
$ws.Range("D42").Value = '1.246'
$ws.Range("E42").Value = '  -0.35%  '

# Row 43
$ws.Range("E43").Value = '  +2.56%  '

# Row 44
$ws.Range("D44").Value = '0.6140'
$ws.Range("E44").Value = '  +1.72%  '

# Row 45
$ws.Range("D45").Value = '13.24'
$ws.Range("E45").Value = '  +1.50%  '

# Row 46
$ws.Range("D46").Value = '3.753'

# Row 47
$ws.Range("D47").Value = '2.058'
$ws.Range("E47").Value = '  +2.49%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '125.88'
$ws.Range("E48").Value = '  +3.17%  '

# Row 49
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.243'
$ws.Range("E49").Value = '  +2.22%  '

# Row 50
$ws.Range("D50").Value = '79.54'
$ws.Range("E50").Value = '  +3.73%  '

# Row 51
$ws.Range("D51").Value = '1.147'
$ws.Range("E51").Value = '  -2.42%  '
